$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").Value = "13`$ID_partida/puntos/numero1-palo1/numero2-palo2"

$ws.Range("F20").Value = "Número palo"
$ws.Range("G20").Value = "Carta"
$ws.Range("F21").Value = "1. Picas`n2. Treboles`n3. Corazones`n4. Rombos"
$ws.Range("F21").WrapText = $true

$ws.Columns("F").ColumnWidth = 13.7109375
$ws.Rows(21).RowHeight = 60

$ws.Range("F22").Select() | Out-Null
